# Weekly update: insert a new price record as a new row 556, pushing the
# existing rows 556:647 down to 557:648 (sheet dimension grows from
# A1:R647 to A1:R648).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 556 (Excel shifts rows 556:647 -> 557:648)
$ws.Rows.Item(556).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A556").Value = 5
$ws.Range("B556").Value = "Macroferia Regional de Talca"
$ws.Range("C556").Value = "Maule"
$ws.Range("D556").Value = 45218
$ws.Range("E556").Value = 7
$ws.Range("F556").Value = 100112023
$ws.Range("G556").Value = "Brócoli"
$ws.Range("H556").Value = "Sin especificar"
$ws.Range("I556").Value = "Primera"
$ws.Range("J556").Value = 4000
$ws.Range("K556").Value = 700
$ws.Range("L556").Value = 700
$ws.Range("M556").Value = 700
$ws.Range("N556").Value = '$/unidad'
$ws.Range("O556").Value = "Región del Maule"
$ws.Range("P556").Value = 700
$ws.Range("Q556").Value = 1
$ws.Range("R556").Value = "Hortaliza"
